$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "Create" + bookmark "_GoBack" + "d By: " -> single run "Created By: "
#    Find.Execute can match text that spans multiple runs (and a
#    bookmark sitting between them); replacing it collapses everything
#    down to one run that keeps the formatting of the first matched
#    run and drops the bookmark entirely - exactly what the diff wants.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Created By: ", $true, $false, $false, $false, $false, $true, 1, $false, "Created By: ", 2) | Out-Null

# ------------------------------------------------------------------
# 2. After the {{description}} paragraph, add:
#      - a Heading2 paragraph "Comments"
#      - a paragraph with {{comments}} split across 3 runs
#      - an empty paragraph that only carries the _GoBack bookmark
# ------------------------------------------------------------------
$descIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq "{{description}}") {
        $descIndex = $i
    }
}

$d.Paragraphs.Item($descIndex).Range.InsertParagraphAfter() | Out-Null
$commentsHeadingIndex = $descIndex + 1
$commentsHeading = $d.Paragraphs.Item($commentsHeadingIndex)
$commentsHeading.Range.Text = "Comments"
$commentsHeading.Range.Style = "Heading 2"

$commentsHeading.Range.InsertParagraphAfter() | Out-Null
$commentsBodyIndex = $commentsHeadingIndex + 1
$commentsBody = $d.Paragraphs.Item($commentsBodyIndex)
$commentsBody.Range.Style = "Normal"
$commentsBody.Range.Text = "{{comments}}"

# Split the single run into three runs ("{{", "comments", "}}") by
# briefly toggling formatting on two sub-ranges - this forces run
# boundaries to appear without altering the visible formatting.
$bodyStart = $commentsBody.Range.Start
$part1 = $d.Range($bodyStart, $bodyStart + 2)
$part1.Bold = 1
$part1.Bold = 0
$part2 = $d.Range($bodyStart + 2, $bodyStart + 10)
$part2.Bold = 1
$part2.Bold = 0

$commentsBody.Range.InsertParagraphAfter() | Out-Null
$bmIndex = $commentsBodyIndex + 1
$bmPara = $d.Paragraphs.Item($bmIndex)
$bmPara.Range.Text = "X"
$bmRange = $d.Range($bmPara.Range.Start, $bmPara.Range.Start + 1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
$bmRange.Text = ""
